# DailyUpdate_Sheet.xlsx — log a new entry: the git commands were copied
# into the embedded interview doc, so add a row for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily_Update")

# Date / Day — same as the rows above (18-07-2015 / Saturday). A brand new
# cell picks up its column's default style automatically (style 2, like the
# other A/B/F/G cells), so a plain value assignment is enough here.
$ws.Cells.Item(4, 1).Value = $ws.Cells.Item(3, 1).Value()
$ws.Cells.Item(4, 2).Value = $ws.Cells.Item(3, 2).Value()

# Start_Time / Stop_Time — 8:00 PM, formatted like the other time cells
# (carry the existing time-formatted style from row 3 across).
$ws.Cells.Item(4, 3).Value = 0.83333333333333337
$ws.Cells.Item(4, 4).Value = 0.83333333333333337
$ws.Range("C3:D3").Copy()
$ws.Range("C4:D4").PasteSpecial(-4122)   # xlPasteFormats

# Learning — wraps across lines like E2, so copy that cell's wrap-text
# formatting over before filling in the text.
$ws.Range("E2").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Cells.Item(4, 5).Value = "copied git commands to`nembedded interview doc"

# SHA1 ID(Previous) / Work
$ws.Cells.Item(4, 6).Value = "fa03bbce022b8909f4d14a627e7aa1915a782d29"
$ws.Cells.Item(4, 7).Value = "check Embedded interview doc"

# The wrapped Learning text spans two lines, so the row grows to fit it.
$ws.Rows.Item(4).RowHeight = 30

# Selection moves on, as it did for the author after the edit.
$ws.Range("G9").Select()
